# Auto-generated edit script applying scheduled market-data refresh to Seraph_Profits sheets.
# Updates currentAveragePrice/NQ/HQ (H/I/J), LevePrice NQ/HQ (K/L), and LeveProfit NQ/HQ (M/N)
# columns for the rows whose underlying market data changed.

$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H74").Value = 9000
$ws.Range("I74").Value = 9000
$ws.Range("K74").Value = 9000
$ws.Range("M74").Value = -8064
$ws.Range("H77").Value = 9000
$ws.Range("I77").Value = 9000
$ws.Range("K77").Value = 45000
$ws.Range("M77").Value = -40320
$ws.Range("H86").Value = 5880.6
$ws.Range("I86").Value = 5701.5
$ws.Range("J86").Value = 6000
$ws.Range("K86").Value = 5701.5
$ws.Range("L86").Value = 6000
$ws.Range("M86").Value = -4578.5
$ws.Range("N86").Value = -8246
$ws.Range("H89").Value = 5880.6
$ws.Range("I89").Value = 5701.5
$ws.Range("J89").Value = 6000
$ws.Range("K89").Value = 28507.5
$ws.Range("L89").Value = 30000
$ws.Range("M89").Value = -22891.5
$ws.Range("N89").Value = -41232
$ws.Range("H116").Value = 7099.8
$ws.Range("I116").Value = 6750
$ws.Range("J116").Value = 7333
$ws.Range("K116").Value = 6750
$ws.Range("L116").Value = 7333
$ws.Range("M116").Value = -3308
$ws.Range("N116").Value = -14217
$ws.Range("H137").Value = 2364.182
$ws.Range("I137").Value = 1976.75
$ws.Range("J137").Value = 2829.1
$ws.Range("K137").Value = 5930.25
$ws.Range("L137").Value = 8487.299999999999
$ws.Range("M137").Value = -3380.25
$ws.Range("N137").Value = -13587.3
$ws.Range("H138").Value = 5765.75
$ws.Range("J138").Value = 5884.3335
$ws.Range("L138").Value = 17653.0005
$ws.Range("N138").Value = -27933.0005

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1936
$ws.Range("I45").Value = 1936
$ws.Range("K45").Value = 1936
$ws.Range("M45").Value = -1559
$ws.Range("H74").Value = 1353.0883
$ws.Range("I74").Value = 892.4838999999999
$ws.Range("K74").Value = 892.4838999999999
$ws.Range("M74").Value = -18.48389999999995
$ws.Range("H77").Value = 1353.0883
$ws.Range("I77").Value = 892.4838999999999
$ws.Range("K77").Value = 4462.4195
$ws.Range("M77").Value = -94.41949999999997
$ws.Range("H92").Value = 76608.8
$ws.Range("I92").Value = 69000
$ws.Range("J92").Value = 78511
$ws.Range("K92").Value = 69000
$ws.Range("L92").Value = 78511
$ws.Range("M92").Value = -66504
$ws.Range("N92").Value = -83503
$ws.Range("H125").Value = 0
$ws.Range("J125").Value = 0
$ws.Range("L125").Value = 0
$ws.Range("N125").ClearContents()
$ws.Range("H132").Value = 2297.7
$ws.Range("I132").Value = 1664.1111
$ws.Range("J132").Value = 8000
$ws.Range("K132").Value = 4992.3333
$ws.Range("L132").Value = 24000
$ws.Range("M132").Value = -2462.3333
$ws.Range("N132").Value = -29060

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 3684
$ws.Range("I105").Value = 3564.3635
$ws.Range("K105").Value = 3564.3635
$ws.Range("M105").Value = -1817.3635
$ws.Range("H134").Value = 3213.4375
$ws.Range("I134").Value = 2994.3635
$ws.Range("J134").Value = 3695.4
$ws.Range("K134").Value = 8983.0905
$ws.Range("L134").Value = 11086.2
$ws.Range("M134").Value = -6448.0905
$ws.Range("N134").Value = -16156.2

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 11695.2
$ws.Range("I99").Value = 7311.9287
$ws.Range("K99").Value = 7311.9287
$ws.Range("M99").Value = -5813.9287
$ws.Range("H126").Value = 11695.2
$ws.Range("I126").Value = 7311.9287
$ws.Range("K126").Value = 21935.7861
$ws.Range("M126").Value = -19465.7861
$ws.Range("H134").Value = 2999.9092
$ws.Range("I134").Value = 2387.2856
$ws.Range("K134").Value = 7161.8568
$ws.Range("M134").Value = -4626.8568

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 66688.2
$ws.Range("I2").Value = 111126.89
$ws.Range("K2").Value = 666761.34
$ws.Range("M2").Value = -666648.34
$ws.Range("H12").Value = 470.72726
$ws.Range("I12").Value = 441.6154
$ws.Range("K12").Value = 1324.8462
$ws.Range("M12").Value = -1151.8462

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 5709.6
$ws.Range("I80").Value = 5099.6665
$ws.Range("K80").Value = 5099.6665
$ws.Range("M80").Value = -4101.6665
$ws.Range("H83").Value = 5709.6
$ws.Range("I83").Value = 5099.6665
$ws.Range("K83").Value = 25498.3325
$ws.Range("M83").Value = -20506.3325
$ws.Range("H126").Value = 4348.5557
$ws.Range("I126").Value = 3280.75
$ws.Range("J126").Value = 5202.8
$ws.Range("K126").Value = 9842.25
$ws.Range("L126").Value = 15608.4
$ws.Range("M126").Value = -7372.25
$ws.Range("N126").Value = -20548.4

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5500
$ws.Range("J7").Value = 6000
$ws.Range("L7").Value = 6000
$ws.Range("N7").Value = -6224
$ws.Range("H68").Value = 2777.8
$ws.Range("I68").Value = 2972.25
$ws.Range("J68").Value = 2000
$ws.Range("K68").Value = 2972.25
$ws.Range("L68").Value = 2000
$ws.Range("M68").Value = -2223.25
$ws.Range("N68").Value = -3498
$ws.Range("H71").Value = 2777.8
$ws.Range("I71").Value = 2972.25
$ws.Range("J71").Value = 2000
$ws.Range("K71").Value = 14861.25
$ws.Range("L71").Value = 10000
$ws.Range("M71").Value = -11117.25
$ws.Range("N71").Value = -17488
$ws.Range("H126").Value = 5500
$ws.Range("J126").Value = 6000
$ws.Range("L126").Value = 18000
$ws.Range("N126").Value = -22940
$ws.Range("H127").Value = 51905
$ws.Range("J127").Value = 51905
$ws.Range("L127").Value = 51905
$ws.Range("N127").Value = -61825
$ws.Range("H132").Value = 3965.2188
$ws.Range("I132").Value = 3489.55
$ws.Range("J132").Value = 4758
$ws.Range("K132").Value = 10468.65
$ws.Range("L132").Value = 14274
$ws.Range("M132").Value = -7938.650000000001
$ws.Range("N132").Value = -19334
$ws.Range("H136").Value = 1685.4286
$ws.Range("I136").Value = 966.5
$ws.Range("J136").Value = 5999
$ws.Range("K136").Value = 2899.5
$ws.Range("L136").Value = 17997
$ws.Range("M136").Value = -349.5
$ws.Range("N136").Value = -23097

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 5000375
$ws.Range("H62").Value = 8399.25
$ws.Range("I62").Value = 8023.75
$ws.Range("K62").Value = 8023.75
$ws.Range("M62").Value = -7399.75
$ws.Range("H65").Value = 8399.25
$ws.Range("I65").Value = 8023.75
$ws.Range("K65").Value = 40118.75
$ws.Range("M65").Value = -36998.75
$ws.Range("H107").Value = 695.4
$ws.Range("I107").Value = 744.25
$ws.Range("J107").Value = 500
$ws.Range("K107").Value = 2232.75
$ws.Range("L107").Value = 1500
$ws.Range("M107").Value = -312.75
$ws.Range("N107").Value = -5340
$ws.Range("H120").Value = 30000
$ws.Range("J120").Value = 0
$ws.Range("L120").Value = 0
$ws.Range("N120").ClearContents()
